$d = $word.ActiveDocument

# Locate the phrase "final year " that precedes the run containing "for"
$rng = $d.Content
$null = $rng.Find.Execute("final year ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$startFinal = $rng.Start          # start of "final year "
$endPhrase  = $rng.End            # end of "final year " (== start of the "for" run)

# "final " is 6 characters, "year" is 4 characters
$posYearStart = $startFinal + 6
$posYearEnd   = $posYearStart + 4

# Insert temporary split-marker bookmarks at the boundaries so that editing the
# "year" sub-range does not re-merge it with neighbouring text runs.
$splitA = $d.Range($posYearStart, $posYearStart)
$d.Bookmarks.Add("ZZZ_SPLIT_A", $splitA)

$splitB = $d.Range($posYearEnd, $posYearEnd)
$d.Bookmarks.Add("ZZZ_SPLIT_B", $splitB)

# Replace "year" with "semester" within the bounded range
$wordRng = $d.Range($posYearStart, $posYearEnd)
$wordRng.Text = "semester"

# Remove the old _GoBack bookmark (it used to sit in the empty paragraph
# near the end of the document); it will be re-created at the new edit location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Recompute the end of "semester" (same start, new length)
$newWordEnd = $posYearStart + 8   # len("semester") == 8

# Insert the _GoBack bookmark immediately after "semester", matching Word's
# behaviour of tracking the most recent edit location.
$goBackRng = $d.Range($newWordEnd, $newWordEnd)
$d.Bookmarks.Add("_GoBack", $goBackRng)

# Remove the temporary helper bookmarks
$d.Bookmarks("ZZZ_SPLIT_A").Delete()
$d.Bookmarks("ZZZ_SPLIT_B").Delete()
